# get named ranges of sheet return cluster of 2d arrays instead of 3d array.
# corrected test herefore: update the "Instruments" overview sheet so the
# identifier values point at the new "#CAN1" / "#PCAN_USBBUS1  0x51" names,
# and make "Instruments" the active/selected sheet (was "CAN1").

$wb = $excel.ActiveWorkbook

$instruments = $wb.Worksheets.Item("Instruments")

# New shared-string values referenced by the Instruments sheet.
$instruments.Range("A2").Value = "#CAN1"
$instruments.Range("A3").Value = "#PCAN_USBBUS1  0x51"

# Make "Instruments" the active sheet/tab (was "CAN1"), with A4 selected.
$instruments.Activate()
$instruments.Range("A4").Select()
